$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, pushing existing rows 101-121 down to 102-122.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with this week's entry.
$ws.Cells.Item(101, 1).Value = 6
$ws.Cells.Item(101, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(101, 3).Value = "Metropolitana"
$ws.Cells.Item(101, 4).Value = 44511
$ws.Cells.Item(101, 5).Value = 13
$ws.Cells.Item(101, 6).Value = 100112029
$ws.Cells.Item(101, 7).Value = "Orégano"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 35
$ws.Cells.Item(101, 11).Value = 8000
$ws.Cells.Item(101, 12).Value = 9000
$ws.Cells.Item(101, 13).Value = 8457
$ws.Cells.Item(101, 14).Value = "$/docena de atados"
$ws.Cells.Item(101, 15).Value = "Región Metropolitana"
$ws.Cells.Item(101, 16).Value = 2819
$ws.Cells.Item(101, 17).Value = 3
$ws.Cells.Item(101, 18).Value = "Hortaliza"
